# Auto-generated Excel COM-interop edit script
# Applies numeric (F/G column) updates across sheets, and a content shift
# in rows 36-46 of the '全部类型' sheet (row 36's event removed, new event appended as row 46).

$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")
$wsShow = $wb.Worksheets.Item("演出")
$wsLocal = $wb.Worksheets.Item("本地生活")
$wsAll = $wb.Worksheets.Item("全部类型")

# ---- 展览 : simple numeric updates ----
$wsExpo.Range("F2").Value = 1298
$wsExpo.Range("F4").Value = 150
$wsExpo.Range("F5").Value = 899
$wsExpo.Range("F9").Value = 830
$wsExpo.Range("F11").Value = 709
$wsExpo.Range("F12").Value = 1377
$wsExpo.Range("F14").Value = 750
$wsExpo.Range("F15").Value = 748
$wsExpo.Range("F19").Value = 648
$wsExpo.Range("F20").Value = 1256
$wsExpo.Range("F24").Value = 5245
$wsExpo.Range("F26").Value = 1012
$wsExpo.Range("F27").Value = 2438
$wsExpo.Range("F28").Value = 5856
$wsExpo.Range("F30").Value = 998
$wsExpo.Range("F31").Value = 594
$wsExpo.Range("F34").Value = 1046
$wsExpo.Range("F38").Value = 684
$wsExpo.Range("F41").Value = 7
$wsExpo.Range("F42").Value = 1082
$wsExpo.Range("F43").Value = 10
$wsExpo.Range("F47").Value = 361
$wsExpo.Range("F49").Value = 20

# ---- 演出 : simple numeric updates ----
$wsShow.Range("F4").Value = 9
$wsShow.Range("F6").Value = 72
$wsShow.Range("F8").Value = 119
$wsShow.Range("F9").Value = 475
$wsShow.Range("F11").Value = 98
$wsShow.Range("F12").Value = 95
$wsShow.Range("F13").Value = 130
$wsShow.Range("F15").Value = 670
$wsShow.Range("G34").Value = 90
$wsShow.Range("F41").Value = 491
$wsShow.Range("F43").Value = 27
$wsShow.Range("F45").Value = 4
$wsShow.Range("F49").Value = 11

# ---- 本地生活 : simple numeric updates ----
$wsLocal.Range("F3").Value = 79
$wsLocal.Range("F6").Value = 393
$wsLocal.Range("F7").Value = 224

# ---- 全部类型 : simple numeric updates ----
$wsAll.Range("F4").Value = 1298
$wsAll.Range("F7").Value = 393
$wsAll.Range("F8").Value = 224
$wsAll.Range("F9").Value = 224
$wsAll.Range("F10").Value = 72
$wsAll.Range("F11").Value = 150
$wsAll.Range("F12").Value = 899
$wsAll.Range("F16").Value = 830
$wsAll.Range("F18").Value = 709
$wsAll.Range("F19").Value = 1377
$wsAll.Range("F20").Value = 98
$wsAll.Range("F21").Value = 95
$wsAll.Range("F23").Value = 750
$wsAll.Range("F24").Value = 130
$wsAll.Range("F25").Value = 748
$wsAll.Range("F26").Value = 1256
$wsAll.Range("F33").Value = 1012
$wsAll.Range("F34").Value = 2438
$wsAll.Range("F35").Value = 5856
$wsAll.Range("F48").Value = 4
$wsAll.Range("F50").Value = 11
# ---- 全部类型 : rows 36-46 content shift ----
# Row 36's original event ("爱乐之城"...) is removed; rows 37-46 each
# move up by one, and a brand-new event is appended as the new row 46.
# NumberFormat is forced to Text on column B first so date-like strings
# (e.g. "2024-05-19") are stored as literal text, matching the source file,
# instead of being auto-parsed into Excel date serials.
$wsAll.Range("B36:B46").NumberFormat = "@"

$wsAll.Range("B36").Value = "2024-05-19"
$wsAll.Range("C36").Value = "上海·中村百合香粉丝见面会"
$wsAll.Range("D36").Value = "茂名南路57号近长乐路 上海兰心大戏院"
$wsAll.Range("E36").Value = "2024.05.19 13:00-05.19 20:30"
$wsAll.Range("F36").Value = 1708
$wsAll.Range("G36").Value = 480
$wsAll.Range("H36").Value = "https://show.bilibili.com/platform/detail.html?id=84235"
$wsAll.Range("I36").Value = "//i2.hdslb.com/bfs/openplatform/202404/jlcJcMk91713169101138.jpeg"

$wsAll.Range("B37").Value = "2024-05-25"
$wsAll.Range("C37").Value = "上海·Anime Market 同人展"
$wsAll.Range("D37").Value = "漕宝路1688号 诺宝中心酒店"
$wsAll.Range("E37").Value = "2024.05.25 08:00-05.26 20:00"
$wsAll.Range("F37").Value = 64
$wsAll.Range("G37").Value = 78
$wsAll.Range("H37").Value = "https://show.bilibili.com/platform/detail.html?id=84737"
$wsAll.Range("I37").Value = "//i1.hdslb.com/bfs/openplatform/202404/DyXYekek1713284815372.png"

$wsAll.Range("B38").Value = "2024-06-01"
$wsAll.Range("C38").Value = "上海·《青城山下·千年等一回》传世国风跨界音乐会"
$wsAll.Range("D38").Value = "南京西路1376号 上海商城剧院"
$wsAll.Range("E38").Value = "2024.06.01 19:30-06.01 21:00"
$wsAll.Range("F38").Value = 4
$wsAll.Range("G38").Value = 90
$wsAll.Range("H38").Value = "https://show.bilibili.com/platform/detail.html?id=83836"
$wsAll.Range("I38").Value = "//i0.hdslb.com/bfs/openplatform/202404/Xgp9kXWX1712125472994.jpeg"

$wsAll.Range("B39").Value = "2024-06-01"
$wsAll.Range("C39").Value = "上海·月遇小马宝莉only"
$wsAll.Range("D39").Value = "淞兴西路248号（复客文创园内） 宝锦宴大酒店"
$wsAll.Range("E39").Value = "2024.06.01 10:00-06.01 15:30"
$wsAll.Range("F39").Value = 31
$wsAll.Range("G39").Value = 19.9
$wsAll.Range("H39").Value = "https://show.bilibili.com/platform/detail.html?id=84693"
$wsAll.Range("I39").Value = "//i1.hdslb.com/bfs/openplatform/202404/EKB0dR0r1713275425895.jpeg"

$wsAll.Range("B40").Value = "2024-06-01"
$wsAll.Range("C40").Value = "上海·第五人格only1.0梦境乐园"
$wsAll.Range("D40").Value = "浦东新区金桥路1599号 上海东方万国企业中心宴会厅"
$wsAll.Range("E40").Value = "2024.06.01 10:00-06.02 17:00"
$wsAll.Range("F40").Value = 684
$wsAll.Range("G40").Value = 95
$wsAll.Range("H40").Value = "https://show.bilibili.com/platform/detail.html?id=83697"
$wsAll.Range("I40").Value = "//i1.hdslb.com/bfs/openplatform/202404/FmgMgZoV1712570114989.jpeg"

$wsAll.Range("B41").Value = "2024-06-08"
$wsAll.Range("C41").Value = "上海·夏日欢愉·羽球节庆·原崩铁同人展"
$wsAll.Range("D41").Value = "鲁班路300号 星光摄影器材城"
$wsAll.Range("E41").Value = "2024.06.08 10:00-06.09 17:00"
$wsAll.Range("F41").Value = 38
$wsAll.Range("G41").Value = 60
$wsAll.Range("H41").Value = "https://show.bilibili.com/platform/detail.html?id=84742"
$wsAll.Range("I41").Value = "//i2.hdslb.com/bfs/openplatform/202404/kbTFe8mJ1713862667234.png"

$wsAll.Range("B42").Value = "2024-06-08"
$wsAll.Range("C42").Value = "上海·第一届妖妖动漫游戏展"
$wsAll.Range("D42").Value = "长宁路1191号来福士西区(W)B1层01号、11号 星零界"
$wsAll.Range("E42").Value = "2024.06.08 10:00-06.10 17:00"
$wsAll.Range("F42").Value = 1082
$wsAll.Range("G42").Value = 68
$wsAll.Range("H42").Value = "https://show.bilibili.com/platform/detail.html?id=84642"
$wsAll.Range("I42").Value = "//i1.hdslb.com/bfs/openplatform/202404/fGytR92V1714112934007.jpeg"

$wsAll.Range("B43").Value = "2024-06-08"
$wsAll.Range("C43").Value = "上海·菊次郎的夏天——久石让钢琴曲梦幻之旅演奏会"
$wsAll.Range("D43").Value = "延安东路523号 凯迪拉克·上海音乐厅"
$wsAll.Range("E43").Value = "2024.06.08 19:30-06.08 21:00"
$wsAll.Range("F43").Value = 39
$wsAll.Range("G43").Value = 80
$wsAll.Range("H43").Value = "https://show.bilibili.com/platform/detail.html?id=81413"
$wsAll.Range("I43").Value = "//i2.hdslb.com/bfs/openplatform/202401/QqKuy4611706169245363.jpeg"

$wsAll.Range("B44").Value = "2024-06-22"
$wsAll.Range("C44").Value = "上海·「多厨狂喜」白金交响乐团二次元交响音乐会"
$wsAll.Range("D44").Value = "丁香路425号 上海东方艺术中心"
$wsAll.Range("E44").Value = "2024.06.22 19:30-06.22 21:30"
$wsAll.Range("F44").Value = 896
$wsAll.Range("G44").Value = 188
$wsAll.Range("H44").Value = "https://show.bilibili.com/platform/detail.html?id=82731"
$wsAll.Range("I44").Value = "//i0.hdslb.com/bfs/openplatform/202403/K3AlF8sr1710230449280.jpeg"

$wsAll.Range("B45").Value = "2024-06-23"
$wsAll.Range("C45").Value = "上海·游马晃祐粉丝见面会"
$wsAll.Range("D45").Value = "宜昌路179号 万代南梦宫上海文化中心"
$wsAll.Range("E45").Value = "2024.06.23 13:00-06.23 20:00"
$wsAll.Range("F45").Value = 491
$wsAll.Range("G45").Value = 480
$wsAll.Range("H45").Value = "https://show.bilibili.com/platform/detail.html?id=84330"
$wsAll.Range("I45").Value = "//i2.hdslb.com/bfs/openplatform/202404/G0vOG4EZ1713257811188.jpeg"

$wsAll.Range("B46").Value = "2024-07-12"
$wsAll.Range("C46").Value = "上海·夜鹿x夜游x真夜中   三夜0nly「夜³歌症候群」联合乐队现场"
$wsAll.Range("D46").Value = "虹许路731号4号楼 THE BOXX•城市乐园"
$wsAll.Range("E46").Value = "2024.07.12 19:30-07.12 23:00"
$wsAll.Range("F46").Value = 27
$wsAll.Range("G46").Value = 99
$wsAll.Range("H46").Value = "https://show.bilibili.com/platform/detail.html?id=85005"
$wsAll.Range("I46").Value = "//i2.hdslb.com/bfs/openplatform/202404/NPObaZdG1714384417870.png"
